# Append the new 2025-12-09 01:22 JST scrape run to the "ランサーズ" sheet.
#
# The existing rows 2..19 are untouched in content (only the "取得日時"
# timestamp in column A advances to the new run time). A brand-new listing
# ("X(旧ツイッター)自動ログインについて", score 13) sorts in right before the
# previous tail of score-13 rows, so it is inserted at row 20 and the three
# rows that used to occupy 20..22 slide down to 21..23. All rows (old and
# new) end up stamped with the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStamp = "2025-12-09 01:22:13"

# 1) Make room for the new listing: push rows 20-22 down to 21-23.
$ws.Range("A20").EntireRow.Insert()

# 2) Populate the newly inserted row with the new listing's data.
$ws.Range("A20").Value = $newStamp
$ws.Range("B20").Value = "X(旧ツイッター)自動ログインについて"
$ws.Range("C20").Value = "システム開発"
$ws.Range("D20").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E20").Value = "期限情報なし"
$ws.Range("F20").Value = "https://www.lancers.jp/work/detail/5449817"
$ws.Range("G20").Value = 13

# 3) Stamp every data row (old + new + shifted) with the new fetch time.
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $newStamp
}

# 4) Rebuild the hyperlinks on column F so each URL cell's link target
#    matches its (possibly shifted) displayed URL, in row order.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
}
